# Apply the "Ventas" report edit described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Sheet1" to "Ventas"
$ws.Name = "Ventas"

# Update header row (row 1)
$ws.Cells.Item(1, 1).Value = "ID"
$ws.Cells.Item(1, 2).Value = "Producto"
$ws.Cells.Item(1, 3).Value = "Ventas"

# Update existing data row 2 (id 1 -> 5, name test1 -> A, value stays 100)
$ws.Cells.Item(2, 1).Value = 5
$ws.Cells.Item(2, 2).Value = "A"
$ws.Cells.Item(2, 3).Value = 100

# Update existing data row 3 (id 2 -> 6, name test2 -> B, value 200 -> 150)
$ws.Cells.Item(3, 1).Value = 6
$ws.Cells.Item(3, 2).Value = "B"
$ws.Cells.Item(3, 3).Value = 150

# Add new data row 4 (id 7, name C, value 200)
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 2).Value = "C"
$ws.Cells.Item(4, 3).Value = 200

Write-Output "edit applied"
